$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the melted "variable" column (L) so each Tweet ID group correctly
# cycles through First Hashtag / Second Hashtag / Third Hashtag instead
# of repeating the same hashtag for every row in the group.
$ws.Range("L5").Value = "Second Hashtag"
$ws.Range("L6").Value = "Third Hashtag"
$ws.Range("L7").Value = "First Hashtag"
$ws.Range("L9").Value = "Third Hashtag"
$ws.Range("L10").Value = "First Hashtag"
$ws.Range("L11").Value = "Second Hashtag"

# L12 no longer keeps its bottom border (matches the plain style used by
# the rest of the L column data rows).
$ws.Range("L12").Borders.Item(9).LineStyle = -4142

# Apply an AutoFilter over the pivoted/melted table, which also records
# the hidden _xlnm._FilterDatabase defined name scoped to Sheet1.
$ws.Range("J3:M12").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$J`$3:`$M`$12")
$filterName.Visible = $false

# Update the view: zoom to 110% and move the active selection/cell.
$excel.ActiveWindow.Zoom = 110
$ws.Range("B10").Select() | Out-Null
